# Get Company Info
# Adds CEO / Media Contact / Auditor / Solicitor / CFO / Website columns
# (C:H) for the first five NZX companies, wraps the long "Media Contact"
# text in column D, and widens columns F/G to fit their contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("C1").Value = "CEO"
$ws.Range("D1").Value = "Media Contact"
$ws.Range("E1").Value = "Auditor"
$ws.Range("F1").Value = "Solicitor"
$ws.Range("G1").Value = "CFO"
$ws.Range("H1").Value = "Website"

# --- Row 2: Ascension Capital Limited ----------------------------------
$ws.Range("D2").Value = "`nJohn Cilliers`nC/o Duncan Cotterill`nLevel 2, Tower Building`n50 Customhouse Quay`nWellington 6143`n+64 9 520 1020`nAscension Capital Limited website"
$ws.Range("E2").Value = "BDO Wellington"
$ws.Range("F2").Value = "Duncan Cotterill"

# --- Row 3: AFC Group Holdings Limited ----------------------------------
$ws.Range("D3").Value = "`nPO Box 230122`nBotany`nAuckland`n+6499300245`nAFC Group Holdings Limited website"
$ws.Range("E3").Value = "William Buck"
$ws.Range("F3").Value = "DLA Piper"
$ws.Range("G3").Value = "Hao Long"
$ws.Range("H3").Value = "http://www.afcnz.com"

# --- Row 4: Australian Foundation Investment Company Limited -----------
$ws.Range("C4").Value = "Mark Freeman"
$ws.Range("D4").Value = "`nAndrew Porter`nMail Box 146,`n101 Collins Street,`nMelbourne, VIC 3000`n+6139650 9911`nAustralian Foundation Investment Company Limited website"
$ws.Range("E4").Value = "PriceWaterhouseCoopers"
$ws.Range("G4").Value = "Andrew Porter"
$ws.Range("H4").Value = "http://www.afi.com.au/"

# --- Row 5: AFT Pharmaceuticals Limited ---------------------------------
$ws.Range("C5").Value = "Hartley Atkinson"
$ws.Range("D5").Value = "`nMalcolm Tubby`nLevel 1`n129 Hurstmere Road`nTakapuna`nAuckland 0622`n+64 9 488 0232`nAFT Pharmaceuticals Limited website"
$ws.Range("E5").Value = "Deloitte"
$ws.Range("F5").Value = "Harmos Horton Lusk Limited"
$ws.Range("G5").Value = "Malcolm Tubby"
$ws.Range("H5").Value = "http://www.aftpharm.com"

# --- Row 6: Smartshares Global Aggregate Bond ETF -----------------------
$ws.Range("D6").Value = "`nJohn McLean`nPO Box 2959 Wellington 6140`n0800808780`nSmartshares Global Aggregate Bond ETF website"
$ws.Range("E6").Value = "KPMG"
$ws.Range("F6").Value = "DLA Piper"
$ws.Range("H6").Value = "https://smartshares.co.nz/"

# --- Formatting ---------------------------------------------------------
# Wrap the long "Media Contact" text so it displays on multiple lines.
$ws.Range("D1:D6").WrapText = $true

# Column widths for the new Solicitor / CFO columns.
$ws.Columns.Item(6).ColumnWidth = 25.42
$ws.Columns.Item(7).ColumnWidth = 26.59

# Row heights so the wrapped text is fully visible.
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 315
$ws.Rows.Item(3).RowHeight = 180
$ws.Rows.Item(4).RowHeight = 330
$ws.Rows.Item(5).RowHeight = 270
$ws.Rows.Item(6).RowHeight = 240

# Restore the selection to "select all" (mirrors Ctrl+A in the source session).
$ws.Range("E9").Select()
$ws.Cells.Select()
